$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.03421546146273613
$ws.Range("C2").Value = 0.01668214239180088
$ws.Range("D2").Value = 0.01089771743863821
$ws.Range("E2").Value = 0.010323453694581985
$ws.Range("F2").Value = 0.00040557500324212015
$ws.Range("J2").Value = 0.12767818570137024
$ws.Range("K2").Value = 1.4655897617340088
